$d = $word.ActiveDocument

# 1) "Use case 3: Hent aktiv produktions data" -- move the "3" into the
#    first run (so the first run reads "Use case 3" and the second run,
#    which carries the da-DK language mark, reads ": Hent aktiv produktions data").
$d.Content.Find.Execute("Use case 3", $true, $false, $false, $false, $false, $true, 1, $false, "Use case 3", 2) | Out-Null

# 2) Same fix for "Use case 4: Send aktiv produktions data".
$d.Content.Find.Execute("Use case 4", $true, $false, $false, $false, $false, $true, 1, $false, "Use case 4", 2) | Out-Null

# 3) Fill in the empty "Reliability" bullet with the new error-handling
#    requirement. The paragraph currently has no runs at all, so we
#    rebuild it (with the same ListParagraph/numPr formatting plus the
#    da-DK language tagging already used by sibling bullets) via
#    Range.InsertXML, which lets us create the two runs exactly as
#    authored (mirrors how Word split the typed sentence into two runs).
$reliabilityPara = $d.Paragraphs(9)
$reliabilityRange = $reliabilityPara.Range
Write-Host "Reliability paragraph before: [$($reliabilityRange.Text)]"

$reliabilityXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>Hvis der ved modtagelse a</w:t></w:r><w:r><w:rPr><w:lang w:val="da-DK"/></w:rPr><w:t>f data er nogen form for konflikt n&#229;r data skal skrives til NAV databasen, skal ingen &#230;ndringer gemmes og en passende fejlmeddelelse returneres.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$reliabilityRange.InsertXML($reliabilityXml)

# 4) The "_GoBack" bookmark (Word's "last edit location" marker) now
#    tracks the Supportability paragraph -- the site of the most recent
#    edit -- instead of the end of the "Log gemt data" paragraph. Move it
#    by deleting the old bookmark and re-creating it, collapsed, right
#    after the "Supportability" text.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$supportPara = $d.Paragraphs(13)
$supportRange = $supportPara.Range
Write-Host "Supportability paragraph before: [$($supportRange.Text)]"

$supportXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Supportability</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$supportRange.InsertXML($supportXml)

Write-Host "Done."
